$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: D1 becomes "median_percent", E1 is removed.
$ws.Range("D1").Value = "median_percent"

# Update data values for C (median_ng_gTEQ) and D (now median_percent).
$ws.Range("C2").Value = 0.00106495
$ws.Range("D2").Value = 16.38384615384615

$ws.Range("C3").Value = 0.0030336
$ws.Range("D3").Value = 46.67076923076923

$ws.Range("C4").Value = 0.017437
$ws.Range("D4").Value = 268.2615384615385

# Remove column E entirely (was median_2), shrinking the used range to A1:D4.
$ws.Range("E1:E4").Delete()
